$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Copy row formatting (bold/centered/bordered id style, date number format)
# from the last existing data row (136) onto the two newly appended rows
# (137, 138) before populating their values, so the style indices for
# column A (id) and column E (Date) match the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A136").Copy() | Out-Null
$ws.Range("A137:A138").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E136").Copy() | Out-Null
$ws.Range("E137:E138").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Swap data for rows 95 / 96 (B:AC only, id in column A unchanged)
$ws.Range("B95").Value = 6482819
$ws.Range("C95").Value = "Estonia Meistriliiga"
$ws.Range("D95").Value = "Estonia Meistriliiga"
$ws.Range("E95").Value = 45231.54166666666
$ws.Range("F95").Value = "JK Tammeka Tartu"
$ws.Range("G95").Value = "FC Kuressaare"
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = "A"
$ws.Range("K95").Value = 1.833
$ws.Range("L95").Value = 3.5
$ws.Range("M95").Value = 3.5
$ws.Range("N95").Value = 2.1
$ws.Range("O95").Value = 3.4
$ws.Range("P95").Value = 2.875
$ws.Range("Q95").Value = -0.25
$ws.Range("R95").Value = 1.975
$ws.Range("S95").Value = 1.825
$ws.Range("T95").Value = 3
$ws.Range("U95").Value = 1.825
$ws.Range("V95").Value = 1.975
$ws.Range("W95").Value = -1
$ws.Range("X95").Value = -1
$ws.Range("Y95").Value = 1.875
$ws.Range("Z95").Value = -1
$ws.Range("AA95").Value = 0.825
$ws.Range("AB95").Value = -1
$ws.Range("AC95").Value = 0.9750000000000001
$ws.Range("B96").Value = 6416370
$ws.Range("C96").Value = "Estonia Meistriliiga"
$ws.Range("D96").Value = "Estonia Meistriliiga"
$ws.Range("E96").Value = 45231.54166666666
$ws.Range("F96").Value = "FC Levadia Tallinn"
$ws.Range("G96").Value = "Parnu JK Vaprus"
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = "D"
$ws.Range("K96").Value = 1.166
$ws.Range("L96").Value = 7
$ws.Range("M96").Value = 11
$ws.Range("N96").Value = 1.2
$ws.Range("O96").Value = 6
$ws.Range("P96").Value = 11
$ws.Range("Q96").Value = -2
$ws.Range("R96").Value = 1.85
$ws.Range("S96").Value = 1.95
$ws.Range("T96").Value = 3
$ws.Range("U96").Value = 1.85
$ws.Range("V96").Value = 1.95
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = 5
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = -1
$ws.Range("AA96").Value = 0.95
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.95

# Swap data for rows 115 / 116 (B:AC only, id in column A unchanged)
$ws.Range("B115").Value = 7919322
$ws.Range("C115").Value = "Estonia Meistriliiga"
$ws.Range("D115").Value = "Estonia Meistriliiga"
$ws.Range("E115").Value = 45360.39583333334
$ws.Range("F115").Value = "FC Kuressaare"
$ws.Range("G115").Value = "FC Levadia Tallinn"
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 6
$ws.Range("J115").Value = "A"
$ws.Range("K115").Value = 11
$ws.Range("L115").Value = 6
$ws.Range("M115").Value = 1.166
$ws.Range("N115").Value = 15
$ws.Range("O115").Value = 8.5
$ws.Range("P115").Value = 1.125
$ws.Range("Q115").Value = 2.5
$ws.Range("R115").Value = 1.825
$ws.Range("S115").Value = 1.975
$ws.Range("T115").Value = 3.25
$ws.Range("U115").Value = 1.9
$ws.Range("V115").Value = 1.9
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = -1
$ws.Range("Y115").Value = 0.125
$ws.Range("Z115").Value = -1
$ws.Range("AA115").Value = 0.9750000000000001
$ws.Range("AB115").Value = 0.8999999999999999
$ws.Range("AC115").Value = -1
$ws.Range("B116").Value = 7919323
$ws.Range("C116").Value = "Estonia Meistriliiga"
$ws.Range("D116").Value = "Estonia Meistriliiga"
$ws.Range("E116").Value = 45360.39583333334
$ws.Range("F116").Value = "JK Nomme Kalju"
$ws.Range("G116").Value = "JK Trans Narva"
$ws.Range("H116").Value = 3
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = "H"
$ws.Range("K116").Value = 1.285
$ws.Range("L116").Value = 5.5
$ws.Range("M116").Value = 6.5
$ws.Range("N116").Value = 1.571
$ws.Range("O116").Value = 4.75
$ws.Range("P116").Value = 4.2
$ws.Range("Q116").Value = -1
$ws.Range("R116").Value = 1.925
$ws.Range("S116").Value = 1.875
$ws.Range("T116").Value = 2.75
$ws.Range("U116").Value = 1.875
$ws.Range("V116").Value = 1.925
$ws.Range("W116").Value = 0.571
$ws.Range("X116").Value = -1
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = 0.925
$ws.Range("AA116").Value = -1
$ws.Range("AB116").Value = 0.4375
$ws.Range("AC116").Value = -0.5

# Swap data for rows 120 / 121 (B:AC only, id in column A unchanged)
$ws.Range("B120").Value = 7721087
$ws.Range("C120").Value = "Estonia Meistriliiga"
$ws.Range("D120").Value = "Estonia Meistriliiga"
$ws.Range("E120").Value = 45367.39583333334
$ws.Range("F120").Value = "Paide Linnameeskond"
$ws.Range("G120").Value = "FC Flora Tallinn"
$ws.Range("H120").Value = 2
$ws.Range("I120").Value = 1
$ws.Range("J120").Value = "H"
$ws.Range("K120").Value = 2.2
$ws.Range("L120").Value = 3.3
$ws.Range("M120").Value = 2.8
$ws.Range("N120").Value = 1.85
$ws.Range("O120").Value = 3.6
$ws.Range("P120").Value = 3.4
$ws.Range("Q120").Value = -0.5
$ws.Range("R120").Value = 1.9
$ws.Range("S120").Value = 1.9
$ws.Range("T120").Value = 2.5
$ws.Range("U120").Value = 1.95
$ws.Range("V120").Value = 1.85
$ws.Range("W120").Value = 0.8500000000000001
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 0.8999999999999999
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = 0.95
$ws.Range("AC120").Value = -1
$ws.Range("B121").Value = 7721007
$ws.Range("C121").Value = "Estonia Meistriliiga"
$ws.Range("D121").Value = "Estonia Meistriliiga"
$ws.Range("E121").Value = 45367.39583333334
$ws.Range("F121").Value = "JK Trans Narva"
$ws.Range("G121").Value = "JK Tammeka Tartu"
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 5
$ws.Range("J121").Value = "A"
$ws.Range("K121").Value = 2.25
$ws.Range("L121").Value = 3.3
$ws.Range("M121").Value = 2.75
$ws.Range("N121").Value = 2.1
$ws.Range("O121").Value = 3.25
$ws.Range("P121").Value = 3
$ws.Range("Q121").Value = -0.25
$ws.Range("R121").Value = 1.875
$ws.Range("S121").Value = 1.925
$ws.Range("T121").Value = 2.5
$ws.Range("U121").Value = 1.825
$ws.Range("V121").Value = 1.975
$ws.Range("W121").Value = -1
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = 2
$ws.Range("Z121").Value = -1
$ws.Range("AA121").Value = 0.925
$ws.Range("AB121").Value = 0.825
$ws.Range("AC121").Value = -1

# New row 137
$ws.Range("A137").Value = 135
$ws.Range("B137").Value = 7721016
$ws.Range("C137").Value = "Estonia Meistriliiga"
$ws.Range("D137").Value = "Estonia Meistriliiga"
$ws.Range("E137").Value = 45396.35416666666
$ws.Range("F137").Value = "Parnu JK Vaprus"
$ws.Range("G137").Value = "FC Flora Tallinn"
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 3
$ws.Range("J137").Value = "A"
$ws.Range("K137").Value = 4.333
$ws.Range("L137").Value = 4
$ws.Range("M137").Value = 1.571
$ws.Range("N137").Value = 5.75
$ws.Range("O137").Value = 4.5
$ws.Range("P137").Value = 1.4
$ws.Range("Q137").Value = 1.25
$ws.Range("R137").Value = 1.9
$ws.Range("S137").Value = 1.9
$ws.Range("T137").Value = 2.75
$ws.Range("U137").Value = 1.85
$ws.Range("V137").Value = 1.95
$ws.Range("W137").Value = -1
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 0.3999999999999999
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = 0.8999999999999999
$ws.Range("AB137").Value = 0.8500000000000001
$ws.Range("AC137").Value = -1

# New row 138
$ws.Range("A138").Value = 136
$ws.Range("B138").Value = 7719650
$ws.Range("C138").Value = "Estonia Meistriliiga"
$ws.Range("D138").Value = "Estonia Meistriliiga"
$ws.Range("E138").Value = 45396.45833333334
$ws.Range("F138").Value = "JK Nomme Kalju"
$ws.Range("G138").Value = "Paide Linnameeskond"
$ws.Range("H138").Value = 2
$ws.Range("I138").Value = 1
$ws.Range("J138").Value = "H"
$ws.Range("K138").Value = 2.375
$ws.Range("L138").Value = 3.4
$ws.Range("M138").Value = 2.5
$ws.Range("N138").Value = 2.45
$ws.Range("O138").Value = 3.6
$ws.Range("P138").Value = 2.375
$ws.Range("Q138").Value = 0
$ws.Range("R138").Value = 1.95
$ws.Range("S138").Value = 1.85
$ws.Range("T138").Value = 2.75
$ws.Range("U138").Value = 1.925
$ws.Range("V138").Value = 1.875
$ws.Range("W138").Value = 1.45
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = -1
$ws.Range("Z138").Value = 0.95
$ws.Range("AA138").Value = -1
$ws.Range("AB138").Value = 0.4625
$ws.Range("AC138").Value = -0.5

# ---------------------------------------------------------------------------
# Dimension (A1:AC138) and shared-string usage counts are maintained
# automatically by the workbook engine when cell values are written.
# ---------------------------------------------------------------------------
